{"js": "// Apply the 7 small text edits described by the diff (stray space\n// characters added/removed by the author while proofreading the\n// requirements list). Each edit is applied via a unique, narrow\n// search string so we never touch unrelated text.\n\nconst edits = [\n  {\n    // Narrow match that stops right after the trailing space following the\n    // period - the document's lone bookmark (Word's internal \"_GoBack\"\n    // marker) sits right after that space, so we avoid folding it into a\n    // larger replaced range.\n    find: \"modify existing bookings made for customers. \",\n    replace: \"modify existing bookings made for customers.\"\n  },\n  {\n    find: \"Staff must be able to view a customers personal details\",\n    replace: \"Staff must be able to view a customerspersonal details\"\n  },\n  {\n    find: \"- Flight manager must be able to modify airport details\",\n    replace: \"- Flight manager must be able to modifyairport details\"\n  },\n  {\n    find: \"- Flight manager must be able to modify route details\",\n    replace: \"- Flight manager must be able to modifyroute details\"\n  },\n  {\n    find: \" -Seating on flights must be categorised into classes\",\n    replace: \"-Seating on flights must be categorised into classes\"\n  },\n  {\n    find: \"- Program should run on windows and linux /\",\n    replace: \"- Program should run on windows and linux/\"\n  },\n  {\n    find: \"includes flights up to a week away\",\n    replace: \"includes flights upto a week away\"\n  }\n];\n\nconst body = context.document.body;\n\nfor (const edit of edits) {\n  const results = body.search(edit.find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(edit.replace, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "# Apply the 7 small text edits described by the diff (stray space\n# characters added/removed by the author while proofreading the\n# requirements list). Each edit is applied via Find/Replace with a\n# unique, narrow search string so we never touch unrelated text.\n\n$d = $word.ActiveDocument\n\n$edits = @(\n    # Narrow match that stops right after the trailing space following the\n    # period - the document's lone bookmark (Word's internal \"_GoBack\"\n    # marker) sits right after that space, and a Find/Replace whose match\n    # range swallows the bookmark would delete it outright.\n    @{ Find = \"modify existing bookings made for customers. \"; Replace = \"modify existing bookings made for customers.\" },\n    @{ Find = \"Staff must be able to view a customers personal details\"; Replace = \"Staff must be able to view a customerspersonal details\" },\n    @{ Find = \"- Flight manager must be able to modify airport details\"; Replace = \"- Flight manager must be able to modifyairport details\" },\n    @{ Find = \"- Flight manager must be able to modify route details\"; Replace = \"- Flight manager must be able to modifyroute details\" },\n    @{ Find = \" -Seating on flights must be categorised into classes\"; Replace = \"-Seating on flights must be categorised into classes\" },\n    @{ Find = \"- Program should run on windows and linux /\"; Replace = \"- Program should run on windows and linux/\" },\n    @{ Find = \"includes flights up to a week away\"; Replace = \"includes flights upto a week away\" }\n)\n\nforeach ($edit in $edits) {\n    $rng = $d.Content\n    $rng.Find.Execute($edit.Find, $false, $false, $false, $false, $false, $true, 1, $false, $edit.Replace, 2)\n}\n"}
